$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8: C8 "9:45AM-4PM" (shift) -> "OFF" ---
# Reuse the existing "OFF" style (same as e.g. D8/E8) then set the text.
$ws.Range("D8").Copy()
$ws.Range("C8").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("C8").Value = "OFF"

# --- Row 12: C12 "9:30AM-4PM" -> "9:45AM-4PM" (format/style unchanged) ---
$ws.Range("C12").Value = "9:45AM-4PM"

# --- Row 15: C15 "OFF" -> "10AM-5PM" (shift) ---
# Reuse the existing shift-time style (same as e.g. B15/D15) then set the text.
$ws.Range("B15").Copy()
$ws.Range("C15").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("C15").Value = "10AM-5PM"

# --- Row 16: C16 "10AM-5PM" -> "OFF" ---
# Reuse the existing "OFF" style (same as e.g. B16/D16) then set the text.
$ws.Range("B16").Copy()
$ws.Range("C16").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("C16").Value = "OFF"

# --- Unassigned Shifts block: move "Bartender, 10AM-4PM" from C27 to new C28 ---
# C27 becomes a duplicate of D27 ("Lifeguard, 9:30AM-4PM"); format (style 7) unchanged.
$ws.Range("C27").Value = "Lifeguard,`n9:30AM-4PM"

# New C28 cell: same style as the rest of row 28, with the text that used to be in C27.
$ws.Range("B28").Copy()
$ws.Range("C28").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("C28").Value = "Bartender,`n10AM-4PM"

$excel.CutCopyMode = $false
